$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows whose "Sending cluster" is "ECs" (old rows 2 and 3).
# This drops "ECs" out of the shared-string pool entirely once no cell
# references it any more, and shifts the former rows 4-7 up to rows 2-5.
$ws.Rows("2:3").Delete()

# Refresh the numeric columns (E:T) of the remaining four rows with the
# newly-computed TPM-derived values. Columns A-D (cluster/gene labels)
# are unchanged by this update.

function Set-RowValues($ws, [int]$row, [double[]]$values) {
    $arr = New-Object 'object[,]' 1, $values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $startCol = 5  # column E
    $endCol = $startCol + $values.Length - 1
    $rng = $ws.Range($ws.Cells.Item($row, $startCol), $ws.Cells.Item($row, $endCol))
    $rng.Value = $arr
}

# Row 2 (was row 4): FAPs -> Efna5/Ephb2 -> FAPs
Set-RowValues $ws 2 @(3, 1, 2.900731333333333, 8.702194, 0.8130494232775288, 0.8130494232775289, 3, 1, 6.346253666666667, 19.038761, 0.921725411846598, 0.9217254118465981, 18.40877686018155, 165.678991741634, 0.7494083145221192, 0.7494083145221194)

# Row 3 (was row 5): FAPs -> Efna5/Ephb2 -> MuSCs
Set-RowValues $ws 3 @(3, 1, 2.900731333333333, 8.702194, 0.8130494232775288, 0.8130494232775289, 3, 1, 0.5389353333333333, 1.616806, 0.07827458815340194, 0.07827458815340194, 1.563306608040444, 14.069759472364, 0.06364110875540953, 0.06364110875540954)

# Row 4 (was row 6): MuSCs -> Efna5/Ephb2 -> FAPs
Set-RowValues $ws 4 @(3, 1, 0.6669870000000001, 2.000961, 0.1869505767224711, 0.1869505767224711, 3, 1, 6.346253666666667, 19.038761, 0.921725411846598, 0.9217254118465981, 4.232868694369, 38.09581824932101, 0.1723170973244787, 0.1723170973244787)

# Row 5 (was row 7): MuSCs -> Efna5/Ephb2 -> MuSCs
Set-RowValues $ws 5 @(3, 1, 0.6669870000000001, 2.000961, 0.1869505767224711, 0.1869505767224711, 3, 1, 0.5389353333333333, 1.616806, 0.07827458815340194, 0.07827458815340194, 0.3594628611740001, 3.235165750566, 0.0146334793979924, 0.0146334793979924)
